$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Total" column header (T1) ---
$ws.Cells.Item(1, 20).Value = "Total"

# --- Add Total column (T) values for existing data rows 2-6 ---
$ws.Cells.Item(2, 20).Value = 89438
$ws.Cells.Item(3, 20).Value = 12949
$ws.Cells.Item(4, 20).Value = 40517
$ws.Cells.Item(5, 20).Value = 15079
$ws.Cells.Item(6, 20).Value = 58955

# --- Row 7: "Outros" category ---
$row7 = @{
    1  = "Outros"
    2  = 6716
    3  = 298
    4  = 372
    5  = 1716
    6  = 2823
    7  = 2643
    8  = 3010
    9  = 3485
    10 = 3802
    11 = 4101
    12 = 4772
    13 = 5369
    14 = 5944
    15 = 6177
    16 = 6415
    17 = 6989
    18 = 24269
    19 = 351
    20 = 89252
}
foreach ($col in $row7.Keys) {
    $ws.Cells.Item(7, $col).Value = $row7[$col]
}

# --- Row 8: "Total" (grand total) row ---
$row8 = @{
    1  = "Total"
    2  = 7463
    3  = 494
    4  = 623
    5  = 2142
    6  = 3599
    7  = 3660
    8  = 4691
    9  = 6167
    10 = 8129
    11 = 10689
    12 = 15166
    13 = 20871
    14 = 26994
    15 = 30547
    16 = 32217
    17 = 33855
    18 = 98417
    19 = 466
    20 = 306190
}
foreach ($col in $row8.Keys) {
    $ws.Cells.Item(8, $col).Value = $row8[$col]
}
